# Update the "dSF" column (F) values to reflect repulled data / mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = 3
$ws.Range("F7").Value = 2
$ws.Range("F8").Value = 1
$ws.Range("F9").Value = 2
$ws.Range("F13").Value = 0
$ws.Range("F19").Value = 4
$ws.Range("F22").Value = -6
$ws.Range("F24").Value = 1
$ws.Range("F25").Value = -2
$ws.Range("F26").Value = -4
$ws.Range("F27").Value = 0
$ws.Range("F28").Value = -8
